$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 104.333336
$ws.Range("I9").Value = 86.318184
$ws.Range("K9").Value = 86.318184
$ws.Range("M9").Value = 82.681816
$ws.Range("H111").Value = 3057.0435
$ws.Range("I111").Value = 2658.5715
$ws.Range("K111").Value = 7975.7145
$ws.Range("M111").Value = -4908.7145
$ws.Range("H125").Value = 6526.5
$ws.Range("J125").Value = 9244.6
$ws.Range("L125").Value = 83201.40000000001
$ws.Range("N125").Value = -88121.40000000001
$ws.Range("H135").Value = 32267.5
$ws.Range("I135").Value = 3238.1428
$ws.Range("J135").Value = 100002.664
$ws.Range("K135").Value = 29143.2852
$ws.Range("L135").Value = 900023.976
$ws.Range("M135").Value = -26608.2852
$ws.Range("N135").Value = -905093.976
$ws.Range("H138").Value = 3052.2092
$ws.Range("J138").Value = 4661.4614
$ws.Range("L138").Value = 13984.3842
$ws.Range("N138").Value = -24264.3842

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1638.1538
$ws.Range("I2").Value = 1754.7273
$ws.Range("J2").Value = 997
$ws.Range("K2").Value = 1754.7273
$ws.Range("L2").Value = 997
$ws.Range("M2").Value = -1641.7273
$ws.Range("N2").Value = -1223
$ws.Range("H45").Value = 2139.5
$ws.Range("I45").Value = 1172.8462
$ws.Range("K45").Value = 1172.8462
$ws.Range("M45").Value = -795.8462
$ws.Range("H102").Value = 2724.25
$ws.Range("I102").Value = 2724.25
$ws.Range("K102").Value = 2724.25
$ws.Range("M102").Value = -1102.25
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H110").Value = 4516
$ws.Range("I110").Value = 5181
$ws.Range("J110").Value = 3407.6667
$ws.Range("K110").Value = 5181
$ws.Range("L110").Value = 3407.6667
$ws.Range("M110").Value = -3136
$ws.Range("N110").Value = -7497.6667
$ws.Range("H116").Value = 1638.1538
$ws.Range("I116").Value = 1754.7273
$ws.Range("J116").Value = 997
$ws.Range("K116").Value = 1754.7273
$ws.Range("L116").Value = 997
$ws.Range("M116").Value = 539.2727
$ws.Range("N116").Value = -5585
$ws.Range("H122").Value = 2756
$ws.Range("I122").Value = 1641
$ws.Range("J122").Value = 3592.25
$ws.Range("K122").Value = 4923
$ws.Range("L122").Value = 10776.75
$ws.Range("M122").Value = -2473
$ws.Range("N122").Value = -15676.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1638.1538
$ws.Range("I3").Value = 1754.7273
$ws.Range("J3").Value = 997
$ws.Range("K3").Value = 1754.7273
$ws.Range("L3").Value = 997
$ws.Range("M3").Value = -1640.7273
$ws.Range("N3").Value = -1225
$ws.Range("H20").Value = 1397.2069
$ws.Range("I20").Value = 1053.05
$ws.Range("K20").Value = 1053.05
$ws.Range("M20").Value = -806.05
$ws.Range("H99").Value = 7000
$ws.Range("I99").Value = 7000
$ws.Range("K99").Value = 7000
$ws.Range("M99").Value = -5502
$ws.Range("H105").Value = 765.5
$ws.Range("I105").Value = 770.6667
$ws.Range("J105").Value = 750
$ws.Range("K105").Value = 770.6667
$ws.Range("L105").Value = 750
$ws.Range("M105").Value = 976.3333
$ws.Range("N105").Value = -4244
$ws.Range("H140").Value = 76998.5
$ws.Range("J140").Value = 76998.5
$ws.Range("L140").Value = 76998.5
$ws.Range("N140").Value = -87358.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3326.9
$ws.Range("I58").Value = 2596.6316
$ws.Range("J58").Value = 4588.273
$ws.Range("K58").Value = 2596.6316
$ws.Range("L58").Value = 4588.273
$ws.Range("M58").Value = -2393.6316
$ws.Range("N58").Value = -4994.273
$ws.Range("H107").Value = 4136.1577
$ws.Range("I107").Value = 4872.9287
$ws.Range("K107").Value = 4872.9287
$ws.Range("M107").Value = -2952.9287
$ws.Range("H132").Value = 1448.4546
$ws.Range("I132").Value = 1448.4546
$ws.Range("K132").Value = 4345.3638
$ws.Range("M132").Value = -1815.3638
$ws.Range("H134").Value = 4107.288
$ws.Range("I134").Value = 3616.5
$ws.Range("K134").Value = 10849.5
$ws.Range("M134").Value = -8314.5
$ws.Range("H136").Value = 3326.9
$ws.Range("I136").Value = 2596.6316
$ws.Range("J136").Value = 4588.273
$ws.Range("K136").Value = 7789.8948
$ws.Range("L136").Value = 13764.819
$ws.Range("M136").Value = -5239.8948
$ws.Range("N136").Value = -18864.819
$ws.Range("H140").Value = 124642.8
$ws.Range("J140").Value = 124642.8
$ws.Range("L140").Value = 124642.8
$ws.Range("N140").Value = -135002.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1748923.8
$ws.Range("I64").Value = 7000
$ws.Range("J64").Value = 2329565
$ws.Range("K64").Value = 21000
$ws.Range("L64").Value = 6988695
$ws.Range("M64").Value = -20730
$ws.Range("N64").Value = -6989235
$ws.Range("H67").Value = 1748923.8
$ws.Range("I67").Value = 7000
$ws.Range("J67").Value = 2329565
$ws.Range("K67").Value = 21000
$ws.Range("L67").Value = 6988695
$ws.Range("M67").Value = -20064
$ws.Range("N67").Value = -6990567

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1308
$ws.Range("I97").Value = 620
$ws.Range("K97").Value = 620
$ws.Range("M97").Value = -124

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1519.5834
$ws.Range("I61").Value = 1519.5834
$ws.Range("K61").Value = 1519.5834
$ws.Range("M61").Value = -1317.5834
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H100").Value = 2166.6667
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959
$ws.Range("H113").Value = 1519.5834
$ws.Range("I113").Value = 1519.5834
$ws.Range("K113").Value = 1519.5834
$ws.Range("M113").Value = 650.4166

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H107").Value = 631.26086
$ws.Range("I107").Value = 594
$ws.Range("K107").Value = 1782
$ws.Range("M107").Value = 138
